$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 13; existing rows 13..55 shift down to 14..56
$ws.Rows.Item(13).Insert()

# Populate the newly inserted row 13 with the new weekly data point
$ws.Cells.Item(13, 1).Value = 8
$ws.Cells.Item(13, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(13, 3).Value = "Coquimbo"
$ws.Cells.Item(13, 4).Value = "2022-07-22"
$ws.Cells.Item(13, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(13, 5).Value = 4
$ws.Cells.Item(13, 6).Value = 100114007
$ws.Cells.Item(13, 7).Value = "Jengibre"
$ws.Cells.Item(13, 8).Value = "Sin especificar"
$ws.Cells.Item(13, 9).Value = "Primera"
$ws.Cells.Item(13, 10).Value = 400
$ws.Cells.Item(13, 11).Value = 15000
$ws.Cells.Item(13, 12).Value = 16000
$ws.Cells.Item(13, 13).Value = 15500
$ws.Cells.Item(13, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(13, 15).Value = "Perú"
$ws.Cells.Item(13, 16).Value = 1192
$ws.Cells.Item(13, 17).Value = 13
$ws.Cells.Item(13, 18).Value = "Hortaliza"
